# Rules!B11 currently shows the shared string "R40" (label of the 4th rule
# row). The commit changes it to the text "1" while keeping the cell's
# existing style (s="23") and shared-string (t="s") type.
#
# A plain  $ws.Range("B11").Value = "1"  assignment would make Excel's
# "smart" numeric detection store it as a *number* (t omitted, <v>1</v>)
# instead of text - not what the diff shows. To force a literal text
# value without disturbing the cell's number format / style, write a
# text-producing formula and then immediately convert it in place to a
# static value via Copy + PasteSpecial(values only).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$target = $ws.Range("B11")

# Use a formula that evaluates to the text "1" so the engine treats the
# result as a string, not a number.
$target.Formula = "=""1"""

# Bake the formula result into a plain value in-place (no leftover
# formula, style/format untouched).
$target.Copy()
$target.PasteSpecial(-4163)  # -4163 = xlPasteValues
